$wb = $excel.ActiveWorkbook

# The data-cleanup sheet ("Sheet2") holds the recode notes that had
# commas used as separators where an "or" reads better / avoids being
# mistaken for a CSV delimiter when the book is exported.
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("B4").Value = "combine variables based on skip pattern (Q3.5 & Q122) " + [char]0x2013 + " can be planned in adv/in discussion w partner/to happen after ideal criteria = YES " + [char]0x2013 + " just happens or left to fate or higher power or natural process meant to be = NO"
$ws.Range("B5").Value = "combine variables based on skip pattern (Q2.2 & Q2.5 for males or Q2.7 & Q2.10 for females)"
$ws.Range("B6").Value = "LowControl = no control or a little; HighControl = complete control or a lot of control"
$ws.Range("B7").Value = "LowControl = no control or a little; HighControl = complete control or a lot of control"

# Leave the cursor where the edits were made, and bring this sheet to
# the front since that's the tab that was active when the file was saved.
$ws.Range("B8").Select()
$ws.Activate()
